# Fruta / hortaliza, semanal
#
# The weekly refresh inserts 4 new observation rows (new rows 453-456,
# holding what used to be the pre-update data of rows 449-452) and then
# overwrites rows 448-452 with corrected/new values. The former last row
# (old row 453) is simply pushed down to row 457, untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert 4 blank rows before the old row 453. Excel shifts
#    the old row 453 down to row 457 and copies formatting (incl. the
#    date number-format on column D) from the row immediately above, so
#    the new rows already carry the right style.
$ws.Rows("453:456").Insert()

# 2) Fill the 4 newly inserted rows with the data that rows 449-452 used
#    to hold prior to this edit.
$newRow453 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44552, 13, 100112030, "Poroto granado", "Sin especificar", "Primera", 600, 35000, 40000, 37917, "`$/saco 25 kilos", "Región Metropolitana", 1517, 25, "Hortaliza")
$newRow454 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44544, 13, 100112030, "Poroto granado", "Sin especificar", "Primera", 400, 40000, 42000, 40850, "`$/malla 25 kilos", "Provincia de Limarí", 1634, 25, "Hortaliza")
$newRow455 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44544, 13, 100112030, "Poroto granado", "Sin especificar", "Primera", 400, 40000, 45000, 42875, "`$/saco 25 kilos", "Región Metropolitana", 1715, 25, "Hortaliza")
$newRow456 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44544, 13, 100112030, "Poroto granado", "Sin especificar", "Primera", 500, 40000, 45000, 43200, "`$/saco 25 kilos", "Región de O'Higgins", 1728, 25, "Hortaliza")

$newRows = @($newRow453, $newRow454, $newRow455, $newRow456)
$startRow = 453
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowData = $newRows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowData[$c]
    }
}

# 3) Update rows 448-452 in place with their new values.
$ws.Range("D448").Value = 44628
$ws.Range("K448").Value = 23000
$ws.Range("L448").Value = 25000
$ws.Range("M448").Value = 24125
$ws.Range("N448").Value = "`$/saco 25 kilos"
$ws.Range("P448").Value = 965

$ws.Range("D449").Value = 44628
$ws.Range("J449").Value = 300
$ws.Range("K449").Value = 23000
$ws.Range("L449").Value = 25000
$ws.Range("M449").Value = 24200
$ws.Range("O449").Value = "Región del Maule"
$ws.Range("P449").Value = 968

$ws.Range("D450").Value = 44628
$ws.Range("I450").Value = "Segunda"
$ws.Range("J450").Value = 250
$ws.Range("K450").Value = 20000
$ws.Range("L450").Value = 20000
$ws.Range("M450").Value = 20000
$ws.Range("N450").Value = "`$/saco 25 kilos"
$ws.Range("O450").Value = "Región Metropolitana"
$ws.Range("P450").Value = 800

$ws.Range("D451").Value = 44628
$ws.Range("I451").Value = "Segunda"
$ws.Range("J451").Value = 100
$ws.Range("K451").Value = 20000
$ws.Range("L451").Value = 20000
$ws.Range("M451").Value = 20000
$ws.Range("O451").Value = "Región del Maule"
$ws.Range("P451").Value = 800

$ws.Range("D452").Value = 44552
$ws.Range("J452").Value = 800
$ws.Range("M452").Value = 42375
$ws.Range("N452").Value = "`$/malla 25 kilos"
$ws.Range("O452").Value = "Región Metropolitana"
$ws.Range("P452").Value = 1695
